# Updates crypto price/volume data to match the latest scrape.
# Rows 25/26 (Monero/Toncoin) also swap position in this refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.017.85"
$ws.Range("E2").Value = "  -2.58%  "

# Row 3
$ws.Range("D3").Value = "1.861.05"
$ws.Range("E3").Value = "  -2.25%  "

# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'306.17"
$ws.Range("E5").Value = "  -2.10%  "

# Row 6
$ws.Range("D6").Value = "'1.0000"
$ws.Range("E6").Value = "  -0.13%  "

# Row 7
$ws.Range("D7").Value = "'0.5106"
$ws.Range("E7").Value = "  +2.47%  "

# Row 8
$ws.Range("D8").Value = "'0.3741"
$ws.Range("E8").Value = "  -0.60%  "

# Row 9
$ws.Range("E9").Value = "  -1.82%  "

# Row 10
$ws.Range("D10").Value = "'0.8882"

# Row 11
$ws.Range("D11").Value = "'20.58"
$ws.Range("E11").Value = "  -2.35%  "

# Row 12
$ws.Range("D12").Value = "'0.07548"
$ws.Range("E12").Value = "  -0.98%  "

# Row 13
$ws.Range("D13").Value = "1.853.05"
$ws.Range("E13").Value = "  -2.61%  "

# Row 14
$ws.Range("D14").Value = "'5.295"
$ws.Range("E14").Value = "  -2.80%  "

# Row 15
$ws.Range("D15").Value = "'89.14"
$ws.Range("E15").Value = "  -2.78%  "

# Row 16
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  -0.11%  "

# Row 17
$ws.Range("D17").Value = "'0.000008371"
$ws.Range("E17").Value = "  -3.76%  "

# Row 18
$ws.Range("D18").Value = "'14.08"
$ws.Range("E18").Value = "  -2.93%  "

# Row 19
$ws.Range("D19").Value = "'0.9997"
$ws.Range("E19").Value = "  -0.09%  "

# Row 20
$ws.Range("D20").Value = "27.078.93"
$ws.Range("E20").Value = "  -2.43%  "

# Row 21
$ws.Range("D21").Value = "'5.055"

# Row 22
$ws.Range("D22").Value = "2.091.87"
$ws.Range("E22").Value = "  -2.79%  "

# Row 23
$ws.Range("D23").Value = "'10.50"
$ws.Range("E23").Value = "  -2.87%  "

# Row 24
$ws.Range("D24").Value = "'6.464"
$ws.Range("E24").Value = "  -1.53%  "

# Row 25
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'1.847"
$ws.Range("E25").Value = "  +0.29%  "

# Row 26
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'149.13"
$ws.Range("E26").Value = "  -2.77%  "

# Row 27
$ws.Range("E27").Value = "  -1.75%  "

# Row 28
$ws.Range("D28").Value = "'2.085"
$ws.Range("E28").Value = "  -5.21%  "

# Row 29
$ws.Range("D29").Value = "'112.83"
$ws.Range("E29").Value = "  -1.64%  "

# Row 30
$ws.Range("E30").Value = "  -3.69%  "

# Row 31
$ws.Range("D31").Value = "'4.657"
$ws.Range("E31").Value = "  -2.58%  "

# Row 32
$ws.Range("D32").Value = "'0.09045"
$ws.Range("E32").Value = "  +1.41%  "

# Row 33
$ws.Range("D33").Value = "'0.05126"
$ws.Range("E33").Value = "  -3.20%  "

# Row 34
$ws.Range("D34").Value = "'3.053"
$ws.Range("E34").Value = "  -3.97%  "

# Row 35
$ws.Range("E35").Value = "  -5.66%  "

# Row 36
$ws.Range("D36").Value = "'0.7318"
$ws.Range("E36").Value = "  -6.23%  "

# Row 37
$ws.Range("D37").Value = "'0.02046"
$ws.Range("E37").Value = "  -1.16%  "

# Row 38
$ws.Range("D38").Value = "'2.502"
$ws.Range("E38").Value = "  -4.41%  "

# Row 39
$ws.Range("D39").Value = "'3.055"
$ws.Range("E39").Value = "  +0.02%  "

# Row 40
$ws.Range("D40").Value = "'1.070"
$ws.Range("E40").Value = "  -1.85%  "

# Row 41
$ws.Range("D41").Value = "'0.5334"
$ws.Range("E41").Value = "  -3.13%  "

# Row 42
$ws.Range("D42").Value = "'6.602"
$ws.Range("E42").Value = "  -2.18%  "

# Row 43
$ws.Range("D43").Value = "'116.75"
$ws.Range("E43").Value = "  +2.24%  "

# Row 44
$ws.Range("D44").Value = "'8.309"
$ws.Range("E44").Value = "  -1.75%  "

# Row 45
$ws.Range("D45").Value = "'0.1471"
$ws.Range("E45").Value = "  -2.48%  "

# Row 46
$ws.Range("D46").Value = "'0.9995"
$ws.Range("E46").Value = "  -0.21%  "

# Row 47
$ws.Range("D47").Value = "'0.4619"
$ws.Range("E47").Value = "  -3.21%  "

# Row 48
$ws.Range("D48").Value = "'10.03"
$ws.Range("E48").Value = "  -4.18%  "

# Row 49
$ws.Range("E49").Value = "  -3.76%  "

# Row 50
$ws.Range("D50").Value = "'36.75"
$ws.Range("E50").Value = "  -0.20%  "

# Row 51
$ws.Range("D51").Value = "'64.04"
$ws.Range("E51").Value = "  -4.35%  "
